$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(19, 8).Value = 2310
$ws.Cells.Item(19, 10).Value = 2310
$ws.Cells.Item(19, 12).Value = 2310
$ws.Cells.Item(19, 14).Value = -2660
$ws.Cells.Item(28, 8).Value = 1558.8572
$ws.Cells.Item(28, 9).Value = 1569.0834
$ws.Cells.Item(28, 10).Value = 1497.5
$ws.Cells.Item(28, 11).Value = 1569.0834
$ws.Cells.Item(28, 12).Value = 1497.5
$ws.Cells.Item(28, 13).Value = -1084.0834
$ws.Cells.Item(28, 14).Value = -2467.5
$ws.Cells.Item(87, 8).Value = 55027.5
$ws.Cells.Item(87, 10).Value = 55027.5
$ws.Cells.Item(87, 12).Value = 55027.5
$ws.Cells.Item(87, 14).Value = -57523.5
$ws.Cells.Item(90, 8).Value = 55027.5
$ws.Cells.Item(90, 10).Value = 55027.5
$ws.Cells.Item(90, 12).Value = 165082.5
$ws.Cells.Item(90, 14).Value = -177562.5
$ws.Cells.Item(137, 8).Value = 2231.2666
$ws.Cells.Item(137, 9).Value = 2092.6365
$ws.Cells.Item(137, 10).Value = 2612.5
$ws.Cells.Item(137, 11).Value = 6277.9095
$ws.Cells.Item(137, 12).Value = 7837.5
$ws.Cells.Item(137, 13).Value = -3727.9095
$ws.Cells.Item(137, 14).Value = -12937.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(55, 8).Value = 32540
$ws.Cells.Item(55, 10).Value = 32540
$ws.Cells.Item(55, 12).Value = 32540
$ws.Cells.Item(55, 14).Value = -33170
$ws.Cells.Item(110, 8).Value = 4407.1875
$ws.Cells.Item(110, 9).Value = 3066.56
$ws.Cells.Item(110, 11).Value = 3066.56
$ws.Cells.Item(110, 13).Value = -1021.56
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3919.0588
$ws.Cells.Item(20, 9).Value = 3500
$ws.Cells.Item(20, 11).Value = 3500
$ws.Cells.Item(20, 13).Value = -3253
$ws.Cells.Item(35, 8).Value = 32678.334
$ws.Cells.Item(35, 10).Value = 34054
$ws.Cells.Item(35, 12).Value = 34054
$ws.Cells.Item(35, 14).Value = -34674
$ws.Cells.Item(82, 8).Value = 13775
$ws.Cells.Item(82, 9).Value = 3933.3333
$ws.Cells.Item(82, 10).Value = 43300
$ws.Cells.Item(82, 11).Value = 3933.3333
$ws.Cells.Item(82, 12).Value = 43300
$ws.Cells.Item(82, 13).Value = -3550.3333
$ws.Cells.Item(82, 14).Value = -44066
$ws.Cells.Item(85, 8).Value = 13775
$ws.Cells.Item(85, 9).Value = 3933.3333
$ws.Cells.Item(85, 10).Value = 43300
$ws.Cells.Item(85, 11).Value = 3933.3333
$ws.Cells.Item(85, 12).Value = 43300
$ws.Cells.Item(85, 13).Value = -2607.3333
$ws.Cells.Item(85, 14).Value = -45952
$ws.Cells.Item(94, 8).Value = 2463.4
$ws.Cells.Item(94, 9).Value = 1765.4615
$ws.Cells.Item(94, 11).Value = 1765.4615
$ws.Cells.Item(94, 13).Value = -1314.4615
$ws.Cells.Item(107, 8).Value = 4162.6
$ws.Cells.Item(107, 9).Value = 4466.6665
$ws.Cells.Item(107, 11).Value = 4466.6665
$ws.Cells.Item(107, 13).Value = -2546.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4505.8945
$ws.Cells.Item(31, 9).Value = 3814.3333
$ws.Cells.Item(31, 10).Value = 5128.3
$ws.Cells.Item(31, 11).Value = 3814.3333
$ws.Cells.Item(31, 12).Value = 5128.3
$ws.Cells.Item(31, 13).Value = -3519.3333
$ws.Cells.Item(31, 14).Value = -5718.3
$ws.Cells.Item(34, 8).Value = 4505.8945
$ws.Cells.Item(34, 9).Value = 3814.3333
$ws.Cells.Item(34, 10).Value = 5128.3
$ws.Cells.Item(34, 11).Value = 3814.3333
$ws.Cells.Item(34, 12).Value = 5128.3
$ws.Cells.Item(34, 13).Value = -3612.3333
$ws.Cells.Item(34, 14).Value = -5532.3
$ws.Cells.Item(41, 8).Value = 24349.5
$ws.Cells.Item(41, 10).Value = 24200
$ws.Cells.Item(41, 12).Value = 24200
$ws.Cells.Item(41, 14).Value = -25056
$ws.Cells.Item(59, 8).Value = 40000
$ws.Cells.Item(59, 10).Value = 40000
$ws.Cells.Item(59, 12).Value = 40000
$ws.Cells.Item(59, 14).Value = -42290
$ws.Cells.Item(60, 8).Value = 14197.8
$ws.Cells.Item(60, 9).Value = 5329.6665
$ws.Cells.Item(60, 10).Value = 27500
$ws.Cells.Item(60, 11).Value = 5329.6665
$ws.Cells.Item(60, 12).Value = 27500
$ws.Cells.Item(60, 13).Value = -4818.6665
$ws.Cells.Item(60, 14).Value = -28522
$ws.Cells.Item(68, 8).Value = 39919.8
$ws.Cells.Item(68, 9).Value = 36000
$ws.Cells.Item(68, 10).Value = 42533
$ws.Cells.Item(68, 11).Value = 36000
$ws.Cells.Item(68, 12).Value = 42533
$ws.Cells.Item(68, 13).Value = -35251
$ws.Cells.Item(68, 14).Value = -44031
$ws.Cells.Item(71, 8).Value = 39919.8
$ws.Cells.Item(71, 9).Value = 36000
$ws.Cells.Item(71, 10).Value = 42533
$ws.Cells.Item(71, 11).Value = 108000
$ws.Cells.Item(71, 12).Value = 127599
$ws.Cells.Item(71, 13).Value = -104256
$ws.Cells.Item(71, 14).Value = -135087
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 376.3846
$ws.Cells.Item(107, 10).Value = 455.7
$ws.Cells.Item(107, 12).Value = 1367.1
$ws.Cells.Item(107, 14).Value = -5207.1
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 51589
$ws.Cells.Item(46, 10).Value = 51589
$ws.Cells.Item(46, 12).Value = 51589
$ws.Cells.Item(46, 14).Value = -51901
$ws.Cells.Item(80, 8).Value = 3783
$ws.Cells.Item(80, 9).Value = 4066.4443
$ws.Cells.Item(80, 11).Value = 4066.4443
$ws.Cells.Item(80, 13).Value = -3068.4443
$ws.Cells.Item(83, 8).Value = 3783
$ws.Cells.Item(83, 9).Value = 4066.4443
$ws.Cells.Item(83, 11).Value = 20332.2215
$ws.Cells.Item(83, 13).Value = -15340.2215
$ws.Cells.Item(92, 8).Value = 14129
$ws.Cells.Item(92, 10).Value = 14129
$ws.Cells.Item(92, 12).Value = 14129
$ws.Cells.Item(92, 14).Value = -17873
$ws.Cells.Item(107, 8).Value = 948
$ws.Cells.Item(107, 9).Value = 951.4286
$ws.Cells.Item(107, 11).Value = 951.4286
$ws.Cells.Item(107, 13).Value = 968.5714
$ws.Cells.Item(132, 8).Value = 1279.6
$ws.Cells.Item(132, 9).Value = 1279.6
$ws.Cells.Item(132, 11).Value = 3838.8
$ws.Cells.Item(132, 13).Value = -1308.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3006.8572
$ws.Cells.Item(7, 9).Value = 3006.8572
$ws.Cells.Item(7, 11).Value = 3006.8572
$ws.Cells.Item(7, 13).Value = -2894.8572
$ws.Cells.Item(16, 8).Value = 1613.1818
$ws.Cells.Item(16, 9).Value = 1510.5555
$ws.Cells.Item(16, 10).Value = 2075
$ws.Cells.Item(16, 11).Value = 1510.5555
$ws.Cells.Item(16, 12).Value = 2075
$ws.Cells.Item(16, 13).Value = -1340.5555
$ws.Cells.Item(16, 14).Value = -2415
$ws.Cells.Item(40, 8).Value = 3624.875
$ws.Cells.Item(40, 9).Value = 3400
$ws.Cells.Item(40, 11).Value = 3400
$ws.Cells.Item(40, 13).Value = -3264
$ws.Cells.Item(46, 8).Value = 18542.3
$ws.Cells.Item(46, 9).Value = 12224.5
$ws.Cells.Item(46, 10).Value = 20121.75
$ws.Cells.Item(46, 11).Value = 12224.5
$ws.Cells.Item(46, 12).Value = 20121.75
$ws.Cells.Item(46, 13).Value = -12036.5
$ws.Cells.Item(46, 14).Value = -20497.75
$ws.Cells.Item(93, 8).Value = 9358.370000000001
$ws.Cells.Item(93, 9).Value = 2058.6
$ws.Cells.Item(93, 11).Value = 2058.6
$ws.Cells.Item(93, 13).Value = -810.5999999999999
$ws.Cells.Item(122, 8).Value = 4832.6665
$ws.Cells.Item(122, 9).Value = 4000
$ws.Cells.Item(122, 10).Value = 4999.2
$ws.Cells.Item(122, 11).Value = 12000
$ws.Cells.Item(122, 12).Value = 14997.6
$ws.Cells.Item(122, 13).Value = -9550
$ws.Cells.Item(122, 14).Value = -19897.6
$ws.Cells.Item(126, 8).Value = 3006.8572
$ws.Cells.Item(126, 9).Value = 3006.8572
$ws.Cells.Item(126, 11).Value = 9020.571599999999
$ws.Cells.Item(126, 13).Value = -6550.571599999999
$ws.Cells.Item(132, 8).Value = 10340
$ws.Cells.Item(132, 9).Value = 10658.097
$ws.Cells.Item(132, 11).Value = 31974.291
$ws.Cells.Item(132, 13).Value = -29444.291
$ws.Cells.Item(136, 8).Value = 3873.5
$ws.Cells.Item(136, 9).Value = 3176.5264
$ws.Cells.Item(136, 11).Value = 9529.5792
$ws.Cells.Item(136, 13).Value = -6979.5792
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 107998.5
$ws.Cells.Item(46, 10).Value = 107998.5
$ws.Cells.Item(46, 12).Value = 107998.5
$ws.Cells.Item(46, 14).Value = -108460.5
$ws.Cells.Item(51, 8).Value = 49999
$ws.Cells.Item(51, 10).Value = 49999
$ws.Cells.Item(51, 12).Value = 49999
$ws.Cells.Item(51, 14).Value = -51019
$ws.Cells.Item(70, 8).Value = 40555.4
$ws.Cells.Item(70, 10).Value = 40555.4
$ws.Cells.Item(70, 12).Value = 40555.4
$ws.Cells.Item(70, 14).Value = -41185.4
$ws.Cells.Item(73, 8).Value = 40555.4
$ws.Cells.Item(73, 10).Value = 40555.4
$ws.Cells.Item(73, 12).Value = 40555.4
$ws.Cells.Item(73, 14).Value = -42739.4
$ws.Cells.Item(132, 8).Value = 3456.7556
$ws.Cells.Item(132, 9).Value = 3835.0322
$ws.Cells.Item(132, 11).Value = 11505.0966
$ws.Cells.Item(132, 13).Value = -8975.096600000001
$ws.Cells.Item(134, 8).Value = 107998.5
$ws.Cells.Item(134, 10).Value = 107998.5
$ws.Cells.Item(134, 12).Value = 323995.5
$ws.Cells.Item(134, 14).Value = -329065.5

Write-Output "applied changes"